$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.020.60"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.21%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.823.40"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.56%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.008"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.34%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.62%  "
$ws.Range("E6").Value = "  -0.27%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4660"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.89%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3660"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.75%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07240"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.91%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8607"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.88%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.86"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.87%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07569"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.34%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.781.69"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.57%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.336"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.97%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.514"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.97%  "
$ws.Range("E16").Value = "  -1.72%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.008"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.20%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008643"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.85%  "
$ws.Range("E19").Value = "  -0.33%  "
$ws.Range("B20").Value = "WrappedBTC"
$ws.Range("C20").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "27.009.67"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.18%  "
$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.50"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.151"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.67%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.54"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.20%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.987.91"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.85%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.59"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.24%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.844"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.82%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.16"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.63%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.056"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.06%  "
$ws.Range("E29").Value = "  -2.55%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "115.34"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.72%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08837"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.76%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.959"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.50%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.427"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.69%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.129"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.11%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7190"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.60%  "
$ws.Range("E36").Value = "  -2.35%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05264"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.46%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01924"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.63%  "
$ws.Range("E39").Value = "  +0.17%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.935"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.47%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.163"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.61%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5163"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.91%  "
$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1627"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.07%  "
$ws.Range("B44").Value = "Frax"
$ws.Range("C44").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8594"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -14.99%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.173"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.63%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4807"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.19%  "
$ws.Range("E47").Value = "  -0.28%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.15"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.86%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "102.83"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.21%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.625"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.81%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06246"
$ws.Range("D51").Style = "Normal"
